$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column T (STORED_PROCEDURE_NAME), shifting
# T:Y -> U:Z, and populate the new header cell with the new config key.
$ws.Columns("T").Insert()
$ws.Range("T1").Value = "STORED_PROCEDURE_SCHEMA"

# Match the width of the new column to its neighbour (column S /
# INDEX_PARTITION_COLUMN), which is what Excel does when a column is
# inserted in the middle of a formatted range.
$ws.Columns("T").ColumnWidth = 30.330729166666668

# Update the view: scroll so column N is the left-most visible column,
# and move the active selection to R6.
$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$win.ScrollRow = 1
$ws.Range("R6").Select() | Out-Null

# Disable multi-threaded/concurrent calculation (calcPr concurrentCalc="0").
$mtc = $excel.MultiThreadedCalculation
$mtc.Enabled = $false
